$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (K2:T2)
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 0.130409
$ws.Range("N2").Value2 = 0.391227
$ws.Range("O2").Value2 = 0.9759813398859937
$ws.Range("P2").Value2 = 0.9759813398859937
$ws.Range("Q2").Value2 = 0.09196520925399999
$ws.Range("R2").Value2 = 0.8276868832859999
$ws.Range("S2").Value2 = 0.9759813398859937
$ws.Range("T2").Value2 = 0.9759813398859937

# Update row 3 values (K3:T3)
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 0.6666666666666666
$ws.Range("M3").Value2 = 0.003209333333333333
$ws.Range("N3").Value2 = 0.009627999999999999
$ws.Range("O3").Value2 = 0.02401866011400631
$ws.Range("P3").Value2 = 0.02401866011400631
$ws.Range("Q3").Value2 = 0.002263241122666666
$ws.Range("R3").Value2 = 0.020369170104
$ws.Range("S3").Value2 = 0.02401866011400631
$ws.Range("T3").Value2 = 0.02401866011400631

# Remove row 4 entirely (data no longer present; also drops the now-unused
# "Resolving-Mac" shared string)
$ws.Rows.Item(4).Delete()
